$wb = $excel.ActiveWorkbook

# --- Sheet "sets": update E4 from 0 to 2 ---
$wsSets = $wb.Worksheets.Item("sets")
$wsSets.Range("E4").Value = 2

# --- Sheet "rallies": append two new rows (88 and 89) ---
$wsRallies = $wb.Worksheets.Item("rallies")

$row88 = @(87, 1, 3, 25, "ADV", "ADVERSÁRIO", 5, "MEIO", "PONTO", "ADV", 24, 1, "0 5 m", "FRENTE", "FRENTE", "FRENTE")
$row89 = @(88, 1, 3, 26, "ADV", "ADVERSÁRIO", 5, "SEGUNDA", "PONTO", "ADV", 24, 2, "0 5 seg", "FRENTE", "FRENTE", "FRENTE")

for ($col = 1; $col -le 16; $col++) {
    $wsRallies.Cells.Item(88, $col).Value = $row88[$col - 1]
    $wsRallies.Cells.Item(89, $col).Value = $row89[$col - 1]
}
